# ---------------------------------------------------------------------------
# Rebuild the worksheet contents/header/style per the target revision.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet" to "Sheet1"
$ws.Name = "Sheet1"

# ---------------------------------------------------------------------------
# Header row (row 1) text values
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Model Name"
$ws.Range("C1").Value = "Exact Precision (Micro Avg)"
$ws.Range("D1").Value = "Exact Recall (Micro Avg)"
$ws.Range("E1").Value = "Exact F1 Score (Micro Avg)"
$ws.Range("F1").Value = "Exact Precision (Macro Avg)"
$ws.Range("G1").Value = "Exact Recall (Macro Avg)"
$ws.Range("H1").Value = "Exact F1 Score (Macro Avg)"
$ws.Range("I1").Value = "Exact Precision (Weighted Avg)"
$ws.Range("J1").Value = "Exact Recall (Weighted Avg)"
$ws.Range("K1").Value = "Exact F1 Score (Weighted Avg)"
$ws.Range("L1").Value = "Partial Precision"
$ws.Range("M1").Value = "Partial Recall"
$ws.Range("N1").Value = "Partial F1 Score"
$ws.Range("O1").Value = "Partial TP"
$ws.Range("P1").Value = "Partial FP"
$ws.Range("Q1").Value = "Partial FN"
$ws.Range("R1").Value = "Support"
$ws.Range("S1").Value = "Accuracy"
$ws.Range("T1").Value = "Result Link"
$ws.Range("U1").Value = "Stats Link"
$ws.Range("V1").Value = "No of GPU Used"
$ws.Range("W1").Value = "Power Consumption"
$ws.Range("X1").Value = "Unnamed: 23"

# ---------------------------------------------------------------------------
# Date-looking text (column A, rows 2-3) must be forced to stay text (not
# auto-converted to a date serial) -- format the cells as Text first, assign
# the values, then put the cell style back to Normal so no stray style is
# left behind in the style table (doing this as one combined A2:A3 range
# operation, rather than per cell, avoids leaving an orphaned style entry).
# ---------------------------------------------------------------------------
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "09/11/2025"
$ws.Range("A3").Value = "09/12/2025"
$ws.Range("A2:A3").Style = "Normal"

# ---------------------------------------------------------------------------
# Row 2 data
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Qwen2.5-32B-Instruct"
$ws.Range("C2").Value = 0.4647887323943662
$ws.Range("D2").Value = 0.3333333333333333
$ws.Range("E2").Value = 0.3882352941176471
$ws.Range("F2").Value = 0.2415822304935208
$ws.Range("G2").Value = 0.1553422051460595
$ws.Range("H2").Value = 0.180737384106598
$ws.Range("I2").Value = 0.5120467495467496
$ws.Range("J2").Value = 0.3333333333333333
$ws.Range("K2").Value = 0.3913137470340021
$ws.Range("L2").Value = 0.5829383886255924
$ws.Range("M2").Value = 0.4155405405405405
$ws.Range("N2").Value = 0.485207100591716
$ws.Range("O2").Value = 123
$ws.Range("P2").Value = 88
$ws.Range("Q2").Value = 173
$ws.Range("R2").Value = 297
$ws.Range("S2").Value = 0.956860395446375
$ws.Range("T2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-32B-Instruct_4_shot.txt"
$ws.Range("U2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-32B-Instruct_4_shot.txt"
$ws.Range("V2").Value = "4 MLGPU"
$ws.Range("W2").Value = "0.117 kWh"
$ws.Range("X2").Value = 3604

# ---------------------------------------------------------------------------
# Row 3 data
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Qwen2.5-32B-Instruct"
$ws.Range("C3").Value = 0.4793388429752066
$ws.Range("D3").Value = 0.3905723905723906
$ws.Range("E3").Value = 0.4304267161410019
$ws.Range("F3").Value = 0.5650976015406163
$ws.Range("G3").Value = 0.3866658053922182
$ws.Range("H3").Value = 0.4380008400060464
$ws.Range("I3").Value = 0.5978224707391374
$ws.Range("J3").Value = 0.3905723905723906
$ws.Range("K3").Value = 0.452173622697443
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.4087837837837838
$ws.Range("N3").Value = 0.449814126394052
$ws.Range("O3").Value = 121
$ws.Range("P3").Value = 121
$ws.Range("Q3").Value = 175
$ws.Range("R3").Value = 297
$ws.Range("S3").Value = 0.9484721390053924
$ws.Range("T3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-32B-Instruct_4_shot.txt"
$ws.Range("U3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-32B-Instruct_4_shot.txt"
$ws.Range("V3").Value = "4 MLGPU"
$ws.Range("W3").Value = "0.099 kWh"

# ---------------------------------------------------------------------------
# Header row formatting: bold font, thin box border all around, center/top
# aligned. Build the combined format once on a scratch cell (so the style
# table only grows by the single final combined style) then copy/paste the
# format onto the whole header range in one shot, finally clearing the
# scratch cell so it doesn't end up part of the used range.
# ---------------------------------------------------------------------------
$helper = $ws.Range("Z1")
$helper.Font.Bold = $true
$helper.Borders.LineStyle = 1
$helper.HorizontalAlignment = -4108
$helper.VerticalAlignment = -4160

$helper.Copy()
$ws.Range("A1:X1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$helper.Clear()

Write-Output "done"
